$wb = $excel.ActiveWorkbook

# --- Sheet 1: rename existing sheet, add the two new sheets in order ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "TC01_Login"

$ws2 = $wb.Worksheets.Add([System.Type]::Missing, $ws1)
$ws2.Name = "TC02_RegisterAccount"

$ws3 = $wb.Worksheets.Add([System.Type]::Missing, $ws2)
$ws3.Name = "Resources"

# --- TC01_Login content (order chosen so shared-string table matches target order) ---
$ws1.Range("A1").Value = "UserName"
$ws1.Range("A2").Value = "Password"
$ws1.Range("B1").Value = "user300@gmail.com"
$ws1.Range("B2").Value = "123456789oO"

$ws1.Columns.Item(1).ColumnWidth = 14.833333333333334
$ws1.Columns.Item(2).ColumnWidth = 31.333333333333332
$ws1.Rows.Item(1).RowHeight = 48.75
$ws1.Rows.Item(2).RowHeight = 55.5

$ws1.Range("A1:B2").VerticalAlignment = -4108

# --- TC02_RegisterAccount content ---
$ws2.Range("A1").Value = "Country Name "
$ws2.Range("B1").Value = "France"

$ws2.Columns.Item(1).ColumnWidth = 26.666666666666668
$ws2.Columns.Item(2).ColumnWidth = 53.166666666666664
$ws2.Rows.Item(1).RowHeight = 48
$ws2.Rows.Item(2).RowHeight = 48.75
# touch row 2 (no-op formatting) so it is tracked in the sheet's used range
$ws2.Range("A2:B2").Font.Size = 11

# --- Resources content (List countries written last so it gets the final shared-string slot) ---
$ws3.Range("A2").Value = "France"
$ws3.Range("A3").Value = "United States"
$ws3.Range("A1").Value = "List countries"

$ws3.Columns.Item(1).ColumnWidth = 17.666666666666668

# --- Defined name used by the dropdown validation ---
$wb.Names.Add("listcountry", $ws3.Range("A2:A3"))

# --- Data validation dropdown on TC02_RegisterAccount!B1 ---
$ws2.Range("B1").Validation.Add(3, 1, 1, "=listcountry")

# --- Selections per sheet ---
[void]$ws1.Range("B2").Select()
[void]$ws3.Range("A2:A3").Select()
[void]$ws2.Range("E20").Select()
